# Add visualization with highcharts in the line files
# This adds new downtime-event rows (Line:8 and Line:3) to the data sheet
# and refreshes the batch "uuid" for existing + new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newUuid = "393f4735-76bf-4e50-807c-e7dbccef7ba7"

# Existing rows 2-5 get the new uuid stamped into column G
$ws.Range("G2:G5").Value = $newUuid

# Column B holds dates formatted as plain text (e.g. "01/09/2024").
# Pre-format the destination cells as Text so Excel doesn't silently
# convert the date-looking strings into date serial numbers.
$ws.Range("B6:B15").NumberFormat = "@"

# New data rows to append (line, tdate, DESCR, st, nd, tgap, uuid, downtime)
$newRows = @(
    @("Line:8 Stage:1", "01/09/2024", "pri cl LA",            45300.42056299769,  45300.42067873842,  0.17, $newUuid, 0.0001157407407407407),
    @("Line:8 Stage:1", "01/09/2024", "pri pH HA",             45300.65476473379,  45300.65488047454,  0.17, $newUuid, 0.0001157407407407407),
    @("Line:3 Stage:1", "12/11/2023", "SEC Cl",                45271.36145825232,  45271.36158556713,  0.18, $newUuid, 0.0001273148148148148),
    @("Line:3 Stage:1", "12/11/2023", "PRI pH",                45271.68537890046,  45271.68549464121,  0.17, $newUuid, 0.0001157407407407407),
    @("Line:3 Stage:1", "12/11/2023", "PRI pH SEC pH",         45271.68549475694,  45271.68769383102,  3.17, $newUuid, 0.002199074074074074),
    @("Line:3 Stage:1", "12/11/2023", "PRI pH&rem SEC pH",     45271.6876965625,   45271.68898128472,  1.85, $newUuid, 0.001284722222222222),
    @("Line:3 Stage:1", "12/11/2023", "PRI pH&rem SEC pH&rem", 45271.68897129629,  45271.68908703703,  0.17, $newUuid, 0.0001157407407407407),
    @("Line:3 Stage:1", "12/11/2023", "PRI pH SEC pH",         45271.6890871875,   45271.68920292824,  0.17, $newUuid, 0.0001157407407407407),
    @("Line:3 Stage:1", "12/11/2023", "PRI pH SEC pH",         45271.89825825232,  45271.8996471412,   2,    $newUuid, 0.001388888888888889),
    @("Line:3 Stage:1", "12/11/2023", "PRI Cl&pH",             45272.09804105324,  45272.09827253472,  0.33, $newUuid, 0.0002314814814814815)
)

$startRow = 6
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]

    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 4).NumberFormat = $ws.Cells.Item(2, 4).NumberFormat

    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 5).NumberFormat = $ws.Cells.Item(2, 5).NumberFormat

    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]

    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 8).NumberFormat = $ws.Cells.Item(2, 8).NumberFormat
}
